$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "service dates added" - the Service Status / Service Date header columns are
# being renamed to Visit Status / Visit Date on the header row (row 3).
$ws.Range("D3").Value = "Visit Status"
$ws.Range("E3").Value = "Visit Date"

# Update the view: scroll down a row and leave the active cell on E8
# (matches the sheetView/selection captured in the saved workbook).
$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
